# Rename the existing sheet and add a second sheet, matching the
# "list-column" / "two-row-header" example workbook used by readxl.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Sheet1" -> "list-column" -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "list-column"

# --- Sheet 2: new "two-row-header" sheet, appended after sheet 1 ------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "two-row-header"

# Header row 1 (variable names)
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "species"
$ws2.Range("C1").Value = "death"
$ws2.Range("D1").Value = "weight"

# Header row 2 (extra annotation), entered in the order that reproduces
# the original shared-string insertion order
$ws2.Range("D2").Value = "(in grams)"
$ws2.Range("B2").Value = "(office supply type)"
$ws2.Range("A2").Value = "(at birth)"
$ws2.Range("C2").Value = "(date is approximate)"

# Data row
$ws2.Range("A3").Value = "Clippy"
$ws2.Range("B3").Value = "paperclip"

# Copy the date cell's format (and value) from Sheet1!B4, which already
# carries the "date of death" number format, so the new cell reuses the
# same style record instead of creating a duplicate one.
$ws1.Range("B4").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws2.Range("C3").Value = 39083

$ws2.Range("D3").Value = 0.9

# --- Selections / active sheet -----------------------------------------
$ws1.Activate()
$ws1.Range("A2:A5").Select()

$ws2.Activate()
$ws2.Range("A1:D1").Select()
